# Generate Report for Handback
# The 1e4c2880-118e-4d5c-8048-b7616a3c6a5b.md file has now been handed back
# successfully (in sync with en-US) for both locales, so the report is
# regenerated: the "Ready for handoff" / stale-error rows flip over to
# "Handed back: in sync with en-US" with fresh handback timestamps, and the
# now-resolved error detail is cleared.

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("E3").Value = "Handed back: in sync with en-US"
$ws.Range("F3").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet ------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("C3").Value = "Handed back: in sync with en-US"
$ws.Range("K3").Value = "2016-08-26 18:50:36"
$ws.Range("P3").Value = ""

# --- de-de sheet ------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("C3").Value = "Handed back: in sync with en-US"
$ws.Range("K3").Value = "2016-08-26 18:50:43"
$ws.Range("P3").Value = ""
